# The commit swaps the embedded theme parts ppt/theme/theme1.xml (used by
# the slide master) and ppt/theme/theme2.xml (used by the notes master):
# theme1.xml goes from the "Integral" palette to the stock "Office" palette,
# and theme2.xml goes from "Office" to "Integral". The font scheme and the
# format scheme (fills/lines/effects) are identical between the two themes,
# so only the 12 theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) actually need to move between the two theme parts.

$p = $ppt.ActivePresentation

# COM RGB long values (0x00BBGGRR) for the "Integral" palette (currently on
# the slide master / theme1.xml).
$integralColors = @(0, 16777215, 5332805, 13754083, 3722137, 3646819, 2412774, 38860, 13611854, 10915127, 2465643, 158642)

# COM RGB long values (0x00BBGGRR) for the stock "Office" palette (currently
# on the notes master / theme2.xml).
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

# Slide master's theme (ppt/theme/theme1.xml): Integral -> Office
$masterScheme = $p.SlideMaster.Theme.ThemeElements.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Item($i).RGB = $officeColors[$i - 1]
}

# Notes master's theme (ppt/theme/theme2.xml): Office -> Integral
$notesScheme = $p.NotesMaster.Theme.ThemeElements.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Item($i).RGB = $integralColors[$i - 1]
}
